$d = $word.ActiveDocument

# Update the date paragraph (wdReplaceOne = 1, scoped to the paragraph range)
$d.Paragraphs.Item(1).Range.Find.Execute("2025-04-25 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-26 Saturday", 1)

# Update table cell contents (row, col are 1-based).
# Use wdReplaceOne (1) rather than wdReplaceAll (2) so that the replace stays
# scoped to the individual cell even when the same text occurs in other cells.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Find.Execute("49÷3=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=10, 7", 1)
$t.Cell(1, 2).Range.Find.Execute("87÷5=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=31, 1", 1)
$t.Cell(1, 3).Range.Find.Execute("10÷8=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "98÷2=49, 0", 1)
$t.Cell(1, 4).Range.Find.Execute("11÷6=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 1)
$t.Cell(1, 5).Range.Find.Execute("15÷2=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "58÷4=14, 2", 1)
$t.Cell(5, 1).Range.Find.Execute("10÷5=2, 0", $true, $false, $false, $false, $false, $true, 1, $false, "18÷4=4, 2", 1)
$t.Cell(5, 2).Range.Find.Execute("33÷5=6, 3", $true, $false, $false, $false, $false, $true, 1, $false, "69÷8=8, 5", 1)
$t.Cell(5, 3).Range.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "22÷6=3, 4", 1)
$t.Cell(5, 4).Range.Find.Execute("87÷9=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "24÷8=3, 0", 1)
$t.Cell(5, 5).Range.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "28÷4=7, 0", 1)
$t.Cell(9, 1).Range.Find.Execute("52÷4=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷2=37, 1", 1)
$t.Cell(9, 2).Range.Find.Execute("75÷5=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=17, 1", 1)
$t.Cell(9, 3).Range.Find.Execute("15÷3=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=11, 1", 1)
$t.Cell(9, 4).Range.Find.Execute("57÷7=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 1)
$t.Cell(9, 5).Range.Find.Execute("19÷9=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "18÷4=4, 2", 1)
$t.Cell(13, 1).Range.Find.Execute("14÷5=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "60÷3=20, 0", 1)
$t.Cell(13, 2).Range.Find.Execute("93÷5=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "49÷9=5, 4", 1)
$t.Cell(13, 3).Range.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=5, 2", 1)
$t.Cell(13, 4).Range.Find.Execute("67÷8=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷7=5, 4", 1)
$t.Cell(13, 5).Range.Find.Execute("92÷9=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "87÷2=43, 1", 1)
$t.Cell(17, 1).Range.Find.Execute("11÷4=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=39, 1", 1)
$t.Cell(17, 2).Range.Find.Execute("92÷9=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "28÷2=14, 0", 1)
$t.Cell(17, 3).Range.Find.Execute("10÷7=1, 3", $true, $false, $false, $false, $false, $true, 1, $false, "37÷3=12, 1", 1)
$t.Cell(17, 4).Range.Find.Execute("94÷8=11, 6", $true, $false, $false, $false, $false, $true, 1, $false, "27÷7=3, 6", 1)
$t.Cell(17, 5).Range.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 1)
